$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$insertAt = 12

# 1) Capture the existing hyperlinks (row number + target URL text) before we
#    touch anything, since inserting a row does not shift the Hyperlinks
#    collection's Range refs automatically in this host.
$existingLinks = @()
for ($i = 1; $i -le $ws.Hyperlinks.Count; $i++) {
    $hl = $ws.Hyperlinks.Item($i)
    $rng = $hl.Range
    $addrStr = [string]$rng.Address()
    $rowPart = $addrStr -replace '^\$B\$', ''
    $rowNum = [int]$rowPart
    $urlText = [string]$rng.Value2
    $existingLinks += , @($rowNum, $urlText)
}

# 2) Insert the new row; this shifts down the cell data/styles of rows 12+.
$ws.Rows.Item($insertAt).Insert()

# 3) Populate the new row's data.
$ws.Range("A" + $insertAt).Value = "Extend time of eviction"
$ws.Range("B" + $insertAt).Value = "https://www.illinoislegalaid.org/legal-information/extend-time-eviction"

# 4) Clear out all (now stale) hyperlinks on the sheet.
$ws.Range("B1").Hyperlinks.Delete()

# 5) Re-add the pre-existing hyperlinks at their shifted row numbers.
foreach ($pair in $existingLinks) {
    $origRow = $pair[0]
    $url = $pair[1]
    $newRow = $origRow
    if ($origRow -ge $insertAt) {
        $newRow = $origRow + 1
    }
    $target = $ws.Cells.Item($newRow, 2)
    $ws.Hyperlinks.Add($target, $url) | Out-Null
}

# 6) Add the hyperlink for the newly inserted row.
$newCell = $ws.Range("B" + $insertAt)
$ws.Hyperlinks.Add($newCell, "https://www.illinoislegalaid.org/legal-information/extend-time-eviction") | Out-Null

# 7) Restore the Hyperlink cell style (Add() can fork a near-duplicate style)
#    and make sure plain (non-hyperlink-object) styled cells keep style "1".
for ($r = 2; $r -le $ws.UsedRange.Rows.Count; $r++) {
    $ws.Cells.Item($r, 2).Style = "Hyperlink"
}

$ws.Range("B10").Select()
